$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "Термін вик" + "о" + "нання" (three separate runs) -> a single run
#    reading "Термін виконання". A literal Find/Replace over the text
#    naturally coalesces the trailing run fragments into the run that
#    receives the replacement text.
# ---------------------------------------------------------------------
$null = $d.Content.Find.Execute("Термін вик", $false, $false, $false, $false, $false, $true, 1, $false, "Термін вик", 2)

# ---------------------------------------------------------------------
# 2) The second "Студент" signature-block cell (gridSpan=3, 2802 dxa /
#    140.1 pt wide) flips its vertical alignment from center to bottom.
#    The first, otherwise-identical cell earlier in the table already
#    uses bottom alignment and must stay untouched, so we scan for the
#    matching cell that is still centered instead of hard-coding a row
#    index.
# ---------------------------------------------------------------------
for ($ti = 1; $ti -le $d.Tables.Count; $ti++) {
    $tbl = $d.Tables.Item($ti)
    for ($r = 1; $r -le $tbl.Rows.Count; $r++) {
        $row = $tbl.Rows.Item($r)
        for ($c = 1; $c -le $row.Cells.Count; $c++) {
            $cell = $row.Cells.Item($c)
            if ([Math]::Abs($cell.Width - 140.1) -lt 0.01) {
                $txt = $cell.Range.Text.TrimEnd([char]13, [char]7)
                if ($txt -eq "Студент" -and $cell.VerticalAlignment -eq 1) {
                    $cell.VerticalAlignment = 3
                }
            }
        }
    }
}
